# Insert a new data row at row 415 (pushing existing rows 415..459 down to 416..460)
# and populate it with the new pineapple price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 415; existing row 415 and everything below shifts down one row.
$ws.Rows.Item(415).Insert()

# Populate the newly inserted row 415 with the new record.
$ws.Cells.Item(415, 1).Value = 5
$ws.Cells.Item(415, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(415, 3).Value = "Maule"
$ws.Cells.Item(415, 4).Value = 45194
$ws.Cells.Item(415, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(415, 5).Value = 7
$ws.Cells.Item(415, 6).Value = "Fruta"
$ws.Cells.Item(415, 7).Value = 100108
$ws.Cells.Item(415, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(415, 9).Value = 100108005
$ws.Cells.Item(415, 10).Value = "Piña"
$ws.Cells.Item(415, 11).Value = "Caramelo"
$ws.Cells.Item(415, 12).Value = "Segunda"
$ws.Cells.Item(415, 13).Value = 200
$ws.Cells.Item(415, 14).Value = 21000
$ws.Cells.Item(415, 15).Value = 21000
$ws.Cells.Item(415, 16).Value = 21000
$ws.Cells.Item(415, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(415, 18).Value = "Ecuador"
$ws.Cells.Item(415, 19).Value = 1500
$ws.Cells.Item(415, 20).Value = 14
